$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) rows 2-37 from serial date 45649 to 45650
for ($row = 2; $row -le 37; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45649) {
        $cell.Value = 45650
    }
}
